# 自动更新Excel文件 - 2025-10-28 23:12:15
# Daily refresh: decrement remaining-days (column E) for every data row.
# When remaining days would drop to 0 (i.e. current value is 1), the cycle
# restarts: E is reset to the total-days value (column D) and the start
# date (column F, stored as yyyymmdd) is rolled forward by that many days.
# Row 36 has a corrupted start date and is left untouched, matching source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 36) {
        continue
    }

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $totalDays = $dCell.Value()
    $remaining = $eCell.Value()

    if ($remaining -gt 1) {
        $eCell.Value = $remaining - 1
    } else {
        $eCell.Value = $totalDays

        $startDate = $fCell.Value()
        $dateStr = [string][int]$startDate
        $parsedDate = [DateTime]::ParseExact($dateStr, "yyyyMMdd", $null)
        $newDate = $parsedDate.AddDays($totalDays)
        $fCell.Value = [int]$newDate.ToString("yyyyMMdd")
    }
}
